$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sectors of Institution, in row order (columns A) as already used in the sheet
$sectors = @(
    "Public, 4-year or above",
    "Private nonprofit, 4-year or above",
    "Private for-profit, 4-year or above",
    "Public, 2-year",
    "Private nonprofit, 2-year",
    "Private for-profit, 2-year",
    "Public, less-than 2-year",
    "Private nonprofit, less-than 2-year",
    "Private for-profit, less-than 2-year"
)

$reportingLocation = "Residence Halls (included in on-campus)"
$offense = "Arrest - Illegal Weapon Possesions"
$dates = @("sum2013", "sum2014", "sum2015")

$row = 29
foreach ($date in $dates) {
    foreach ($sector in $sectors) {
        $ws.Cells.Item($row, 1).Value = "'" + $sector
        $ws.Cells.Item($row, 2).Value = "'" + $reportingLocation
        $ws.Cells.Item($row, 3).Value = $offense
        $ws.Cells.Item($row, 4).Value = "'" + $date
        $row = $row + 1
    }
}

$ws.Range("C39").Select()
